$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20:150 down to 21:151
$ws.Rows("20").Insert()

# Populate the new row 20 with the new daily price record
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44635
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112003
$ws.Cells.Item(20, 7).Value = "Ajo"
$ws.Cells.Item(20, 8).Value = "Chino"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 270
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 17000
$ws.Cells.Item(20, 13).Value = 16556
$ws.Cells.Item(20, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(20, 15).Value = "China"
$ws.Cells.Item(20, 16).Value = 1656
$ws.Cells.Item(20, 17).Value = 10
$ws.Cells.Item(20, 18).Value = "Hortaliza"
